$wb = $excel.ActiveWorkbook

# --- Agency sheet: update jurisdiction/agency names, add two new rows ---
$ws = $wb.Worksheets.Item("Agency")
$ws.Activate()
$ws.Range("B2").Value = "Brookfield County Sheriff"
$ws.Range("B3").Value = "Logan City PD"
$ws.Range("B4").Value = "Scriba City PD"
$ws.Range("B5").Value = "Edwards City PD"
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Allentown Parks/Recreation Dept"
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Morgan College Campus Police"
$ws.Range("A2:B7").Select()

# --- AssessmentCategoryType sheet: update category names, add one new row ---
$ws = $wb.Worksheets.Item("AssessmentCategoryType")
$ws.Activate()
$ws.Range("B3").Value = "Danger To Self"
$ws.Range("B4").Value = "Danger To Others"
$ws.Range("B5").Value = "Gravely Disabled"
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Persistently and Acutely Disabled"
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Other"
$ws.Range("B11").Select()

# --- JurisdictionType sheet: update jurisdiction names, remove last row ---
$ws = $wb.Worksheets.Item("JurisdictionType")
$ws.Activate()
$ws.Range("B2").Value = "Brookfield County Superior Court"
$ws.Range("B3").Value = "Logan Municipal Court"
$ws.Range("B4").Value = "Scriba Municipal Court"
$ws.Rows.Item(5).Delete()
$ws.Range("B2:B4").Select()

# Leave AssessmentCategoryType as the active/selected sheet
$wb.Worksheets.Item("AssessmentCategoryType").Activate()
